$d = $word.ActiveDocument

function Get-ParaByText($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { throw "Text not found: $text" }
    $target = $rng.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $target -and $target -lt $p.Range.End) {
            return $p
        }
    }
    throw "Paragraph not found for text: $text"
}

function Delete-ParasBetween($firstText, $lastText) {
    $pStart = Get-ParaByText($firstText)
    $pEnd = Get-ParaByText($lastText)
    $delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $delRange.Delete()
}

function Insert-ParaAfter($anchorPara, $text, $styleName) {
    $idx = $anchorPara.Index + 1
    $anchorPara.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Text = $text
    if ($styleName) { $p.Range.Style = $styleName } else { $p.Range.Style = "Normal" }
    return $p
}

# DELETE before[2:2]
Delete-ParasBetween '+1 (512) 555-0123 | dheeraj@dheerajchand.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX' '+1 (512) 555-0123 | dheeraj@dheerajchand.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX'

# REPLACE before[4:4] with after[3:3]
$p = Get-ParaByText('Experienced data scientist and software engineer with 15+ years of expertise in geospatial analysis, demographic research, and political data. Proven track record of building scalable systems, conducting complex analyses, and delivering actionable insights for campaigns, organizations, and government agencies.')
$p.Range.Text = 'Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide.'

# REPLACE before[6:6] with after[5:5]
$p = Get-ParaByText('CODE • COMPUTE • INTERACT • MEASURE • PLATFORMS • TRACK')
$p.Range.Text = ''

# REPLACE before[8:13] with after[7:42]
Delete-ParasBetween 'Partner - Siege Analytics (Austin, TX) | 2020 - Present' 'Senior Data Scientist - Lake Research Partners (Washington, DC) | 2018 - 2020'
$anchor = Get-ParaByText('PROFESSIONAL EXPERIENCE')
$anchor = Insert-ParaAfter $anchor 'Partner - Siege Analytics (Austin, TX) | 2005 - Present' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Data, Technology and Strategy Consulting' $null
$anchor = Insert-ParaAfter $anchor '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%' $null
$anchor = Insert-ParaAfter $anchor '• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration' $null
$anchor = Insert-ParaAfter $anchor '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%' $null
$anchor = Insert-ParaAfter $anchor 'Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Civic Graph & Civic Pulse Director' $null
$anchor = Insert-ParaAfter $anchor '• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics' $null
$anchor = Insert-ParaAfter $anchor '• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions' $null
$anchor = Insert-ParaAfter $anchor '• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture' $null
$anchor = Insert-ParaAfter $anchor 'Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Big Data Engineering Transformation' $null
$anchor = Insert-ParaAfter $anchor '• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS' $null
$anchor = Insert-ParaAfter $anchor '• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed' $null
$anchor = Insert-ParaAfter $anchor '• Rewrote mission and offerings of department and drafted integration plan with strategy team' $null
$anchor = Insert-ParaAfter $anchor 'Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'SimCrisis Product Owner/Engineer' $null
$anchor = Insert-ParaAfter $anchor '• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement' $null
$anchor = Insert-ParaAfter $anchor '• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies' $null
$anchor = Insert-ParaAfter $anchor '• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures' $null
$anchor = Insert-ParaAfter $anchor 'Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'RACSO Product Owner/Engineer' $null
$anchor = Insert-ParaAfter $anchor '• Designed comprehensive survey instruments for specialized voting segments and niche markets' $null
$anchor = Insert-ParaAfter $anchor '• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis' $null
$anchor = Insert-ParaAfter $anchor '• Wrote RFP and analyzed bids from 1,200 vendors for research platform development' $null
$anchor = Insert-ParaAfter $anchor 'Research Director - PCCC (Washington, DC) | 2010 - 2012' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Political Research & Data Analysis (FLEEM System)' $null
$anchor = Insert-ParaAfter $anchor '• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys' $null
$anchor = Insert-ParaAfter $anchor '• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren' $null
$anchor = Insert-ParaAfter $anchor '• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver' $null
$anchor = Insert-ParaAfter $anchor 'Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Geospatial CRM Development' $null
$anchor = Insert-ParaAfter $anchor '• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously' $null
$anchor = Insert-ParaAfter $anchor '• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers' $null
$anchor = Insert-ParaAfter $anchor '• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill' $null
$anchor = Insert-ParaAfter $anchor 'Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008' 'Heading3'

# REPLACE before[15:29] with after[44:46]
Delete-ParasBetween '• Trained staff on building Python tooling for report generation and analysis' '• Trained staff on PHP/MySQL for data analysis and reporting systems'
$anchor = Get-ParaByText('Political Research & Analytics')
$anchor = Insert-ParaAfter $anchor '• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party' $null
$anchor = Insert-ParaAfter $anchor '• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems' $null
$anchor = Insert-ParaAfter $anchor '• Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+' $null

# REPLACE before[31:34] with after[48:59]
Delete-ParasBetween 'Polling Consortium Dataset Meta-Analysis (2013 - 2016)' 'Impact: Created $400M dataset that became foundation for modern electoral analytics, estimated current value exceeds $1B'
$anchor = Get-ParaByText('KEY PROJECTS')
$anchor = Insert-ParaAfter $anchor 'National Redistricting Platform (2020 - 2021)' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Cloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide' $null
$anchor = Insert-ParaAfter $anchor 'Technologies: GeoDjango, PostGIS, AWS, Docker, React, Python' $null
$anchor = Insert-ParaAfter $anchor 'Impact: Reduced mapping costs by 73.5%, saving organizations $4.7M in operational expenses' $null
$anchor = Insert-ParaAfter $anchor 'FLEEM Political Polling System (2010 - 2012)' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Completely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity' $null
$anchor = Insert-ParaAfter $anchor 'Technologies: Twilio API, Python, Django, PostgreSQL, JavaScript' $null
$anchor = Insert-ParaAfter $anchor 'Impact: Saved $840K in operational costs plus millions in avoided software licensing' $null
$anchor = Insert-ParaAfter $anchor 'Geospatial Demographic Classification System (2013 - 2016)' 'Heading3'
$anchor = Insert-ParaAfter $anchor 'Machine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%' $null
$anchor = Insert-ParaAfter $anchor 'Technologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow' $null
$anchor = Insert-ParaAfter $anchor 'Impact: Corrected demographic data affecting all Black and Asian-American voters nationwide' $null

# REPLACE before[36:44] with after[61:65]
Delete-ParasBetween 'Data Discovery' '• Interfaced with Government and Activism APIs for seamless data integration'
$anchor = Get-ParaByText('KEY ACHIEVEMENTS AND IMPACT')
$anchor = Insert-ParaAfter $anchor 'Impact' 'Heading3'
$anchor = Insert-ParaAfter $anchor '• Discovered systematic race coding errors affecting all Black and Asian-American voters' $null
$anchor = Insert-ParaAfter $anchor '• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M' $null
$anchor = Insert-ParaAfter $anchor '• Built redistricting platform used by thousands of analysts nationwide' $null
$anchor = Insert-ParaAfter $anchor '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%' $null

# DELETE before[46:52]
Delete-ParasBetween 'CODE Python; R; SQL; JavaScript; PHP' 'For a more detailed, full description of my experience, please visit my LinkedIn and Personal Site.'
